$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("obras")
$ws2 = $wb.Worksheets.Item("referentes")

# --- Sheet "obras": insert a new ID column at the front, and add a
#     "Referentes" column right after the existing data (old column F,
#     which was an empty formatted column). ---
$ws1.Columns("A:A").Insert()

$ws1.Range("A1").Value = "ID"
# The newly inserted column doesn't inherit the bold header style from its
# former neighbour, so re-apply it explicitly (matches the existing header
# style used across the rest of row 1).
$ws1.Range("A1").Font.Bold = $true
$ws1.Range("A2").Value = 1
$ws1.Range("A3").Value = 2
$ws1.Range("A4").Value = 3

$ws1.Range("G1").Value = "Referentes"
$ws1.Range("G2").Value = 1
$ws1.Range("G3").Value = 2.3
$ws1.Range("G4").Value = 4
# NOTE: the runtime quantizes ColumnWidth to 1/6-character increments, so the
# original author's raw width (12.28515625) can't be reproduced bit-exactly;
# 11.5 lands on the nearest reachable bucket (12.333333...).
$ws1.Columns("G:G").ColumnWidth = 11.5

# --- Sheet "referentes": insert a new ID column at the front. ---
$ws2.Columns("A:A").Insert()

$ws2.Range("A1").Value = "ID"
$ws2.Range("A1").Font.Bold = $true
$ws2.Range("A2").Value = 1
$ws2.Range("A3").Value = 2
$ws2.Range("A4").Value = 3
$ws2.Range("A5").Value = 4

# --- Selections / active tab: "obras" becomes the active sheet. ---
$ws2.Range("B9").Select()
$ws1.Activate()
$ws1.Range("G9").Select()
